# Scheduled-runner update: refresh market-price-derived columns (H:N)
# on the per-sheet Leve tables. Values below are the latest pulled
# currentAveragePrice* / LevePrice* / LeveProfit* figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 821.86664
$ws.Range("I43").Value = 1008.5
$ws.Range("J43").Value = 697.44446
$ws.Range("K43").Value = 1008.5
$ws.Range("L43").Value = 697.44446
$ws.Range("M43").Value = -939.5
$ws.Range("N43").Value = -835.44446

$ws.Range("H53").Value = 569.36
$ws.Range("I53").Value = 162.23529
$ws.Range("J53").Value = 1434.5
$ws.Range("K53").Value = 162.23529
$ws.Range("L53").Value = 1434.5
$ws.Range("M53").Value = 474.76471
$ws.Range("N53").Value = -2708.5

$ws.Range("H62").Value = 3474.524
$ws.Range("I62").Value = 3673.0625
$ws.Range("J62").Value = 2839.2
$ws.Range("K62").Value = 3673.0625
$ws.Range("L62").Value = 2839.2
$ws.Range("M62").Value = -3049.0625
$ws.Range("N62").Value = -4087.2

$ws.Range("H65").Value = 3474.524
$ws.Range("I65").Value = 3673.0625
$ws.Range("J65").Value = 2839.2
$ws.Range("K65").Value = 18365.3125
$ws.Range("L65").Value = 14196
$ws.Range("M65").Value = -15245.3125
$ws.Range("N65").Value = -20436

$ws.Range("H100").Value = 3163.3333
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3163.3333
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 3163.3333
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -4245.3333

$ws.Range("H101").Value = 293.33334
$ws.Range("I101").Value = 200
$ws.Range("J101").Value = 386.66666
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 1159.99998
$ws.Range("M101").Value = 1022
$ws.Range("N101").Value = -4403.999980000001

$ws.Range("H113").Value = 3966.9048
$ws.Range("I113").Value = 2895.625
$ws.Range("J113").Value = 4626.154
$ws.Range("K113").Value = 2895.625
$ws.Range("L113").Value = 4626.154
$ws.Range("M113").Value = 358.375
$ws.Range("N113").Value = -11134.154

$ws.Range("H129").Value = 1416.6285
$ws.Range("I129").Value = 780.5
$ws.Range("J129").Value = 1605.1111
$ws.Range("K129").Value = 2341.5
$ws.Range("L129").Value = 4815.3333
$ws.Range("M129").Value = 2658.5
$ws.Range("N129").Value = -14815.3333

$ws.Range("H137").Value = 7693889.5
$ws.Range("I137").Value = 12822229
$ws.Range("J137").Value = 1380.1923
$ws.Range("K137").Value = 38466687
$ws.Range("L137").Value = 4140.5769
$ws.Range("M137").Value = -38464137
$ws.Range("N137").Value = -9240.5769

$ws.Range("H138").Value = 2127.2173
$ws.Range("J138").Value = 2136.6042
$ws.Range("L138").Value = 6409.812600000001
$ws.Range("N138").Value = -16689.8126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3273.8
$ws.Range("I2").Value = 1358.8572
$ws.Range("J2").Value = 7742
$ws.Range("K2").Value = 1358.8572
$ws.Range("L2").Value = 7742
$ws.Range("M2").Value = -1245.8572
$ws.Range("N2").Value = -7968

$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 5000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -5350

$ws.Range("H32").Value = 10441838
$ws.Range("I32").Value = 15383.76
$ws.Range("J32").Value = 47679176
$ws.Range("K32").Value = 15383.76
$ws.Range("L32").Value = 47679176
$ws.Range("M32").Value = -15096.76
$ws.Range("N32").Value = -47679750

$ws.Range("H45").Value = 4168.4
$ws.Range("I45").Value = 3289.75
$ws.Range("J45").Value = 4487.909
$ws.Range("K45").Value = 3289.75
$ws.Range("L45").Value = 4487.909
$ws.Range("M45").Value = -2912.75
$ws.Range("N45").Value = -5241.909

$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

$ws.Range("H116").Value = 3273.8
$ws.Range("I116").Value = 1358.8572
$ws.Range("J116").Value = 7742
$ws.Range("K116").Value = 1358.8572
$ws.Range("L116").Value = 7742
$ws.Range("M116").Value = 935.1428000000001
$ws.Range("N116").Value = -12330

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3273.8
$ws.Range("I3").Value = 1358.8572
$ws.Range("J3").Value = 7742
$ws.Range("K3").Value = 1358.8572
$ws.Range("L3").Value = 7742
$ws.Range("M3").Value = -1244.8572
$ws.Range("N3").Value = -7970

$ws.Range("H24").Value = 366.57144
$ws.Range("I24").Value = 366.57144
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 366.57144
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -131.57144
$ws.Range("N24").ClearContents()

$ws.Range("H107").Value = 2212.9412
$ws.Range("I107").Value = 2163.75
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2163.75
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -243.75
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4763801
$ws.Range("I31").Value = 1525.3143
$ws.Range("J31").Value = 9526077
$ws.Range("K31").Value = 1525.3143
$ws.Range("L31").Value = 9526077
$ws.Range("M31").Value = -1230.3143
$ws.Range("N31").Value = -9526667

$ws.Range("H34").Value = 4763801
$ws.Range("I34").Value = 1525.3143
$ws.Range("J34").Value = 9526077
$ws.Range("K34").Value = 1525.3143
$ws.Range("L34").Value = 9526077
$ws.Range("M34").Value = -1323.3143
$ws.Range("N34").Value = -9526481

$ws.Range("H141").Value = 33749
$ws.Range("J141").Value = 33749
$ws.Range("L141").Value = 33749
$ws.Range("N141").Value = -44109

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 555401.2
$ws.Range("I68").Value = 766.9706
$ws.Range("J68").Value = 948267.1
$ws.Range("K68").Value = 2300.9118
$ws.Range("L68").Value = 2844801.3
$ws.Range("M68").Value = -1489.9118
$ws.Range("N68").Value = -2846423.3

$ws.Range("H71").Value = 555401.2
$ws.Range("I71").Value = 766.9706
$ws.Range("J71").Value = 948267.1
$ws.Range("K71").Value = 6902.7354
$ws.Range("L71").Value = 8534403.9
$ws.Range("M71").Value = -2846.7354
$ws.Range("N71").Value = -8542515.9

$ws.Range("H113").Value = 489.56668
$ws.Range("I113").Value = 472.2857
$ws.Range("J113").Value = 504.6875
$ws.Range("K113").Value = 1416.8571
$ws.Range("L113").Value = 1514.0625
$ws.Range("M113").Value = 753.1428999999998
$ws.Range("N113").Value = -5854.0625

$ws.Range("H125").Value = 6065.517
$ws.Range("J125").Value = 6330.769
$ws.Range("L125").Value = 18992.307
$ws.Range("N125").Value = -28832.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 7000
$ws.Range("J17").Value = 7000
$ws.Range("L17").Value = 7000
$ws.Range("N17").Value = -7336

$ws.Range("H20").Value = 67338.336
$ws.Range("J20").Value = 67338.336
$ws.Range("L20").Value = 67338.336
$ws.Range("N20").Value = -67828.336

$ws.Range("H24").Value = 20000000
$ws.Range("I24").Value = 20000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 20000000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -19999827
$ws.Range("N24").ClearContents()

$ws.Range("H126").Value = 17859378
$ws.Range("I126").Value = 41667520
$ws.Range("J126").Value = 3271.75
$ws.Range("K126").Value = 125002560
$ws.Range("L126").Value = 9815.25
$ws.Range("M126").Value = -125000090
$ws.Range("N126").Value = -14755.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1147.4
$ws.Range("I22").Value = 1198.5714
$ws.Range("J22").Value = 1102.625
$ws.Range("K22").Value = 1198.5714
$ws.Range("L22").Value = 1102.625
$ws.Range("M22").Value = -903.5714
$ws.Range("N22").Value = -1692.625

$ws.Range("H27").Value = 1147.4
$ws.Range("I27").Value = 1198.5714
$ws.Range("J27").Value = 1102.625
$ws.Range("K27").Value = 1198.5714
$ws.Range("L27").Value = 1102.625
$ws.Range("M27").Value = -1091.5714
$ws.Range("N27").Value = -1316.625

$ws.Range("H46").Value = 527214.8
$ws.Range("I46").Value = 692.2857
$ws.Range("J46").Value = 2001477.8
$ws.Range("K46").Value = 692.2857
$ws.Range("L46").Value = 2001477.8
$ws.Range("M46").Value = -504.2857
$ws.Range("N46").Value = -2001853.8

$ws.Range("H55").Value = 313.6216
$ws.Range("I55").Value = 308.25
$ws.Range("J55").Value = 319.94116
$ws.Range("K55").Value = 308.25
$ws.Range("L55").Value = 319.94116
$ws.Range("M55").Value = -135.25
$ws.Range("N55").Value = -665.9411600000001

$ws.Range("H122").Value = 35500
$ws.Range("I122").Value = 55000
$ws.Range("J122").Value = 16000
$ws.Range("K122").Value = 165000
$ws.Range("L122").Value = 48000
$ws.Range("M122").Value = -162550
$ws.Range("N122").Value = -52900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 102506.75
$ws.Range("I15").Value = 200006
$ws.Range("J15").Value = 70007
$ws.Range("K15").Value = 200006
$ws.Range("L15").Value = 70007
$ws.Range("M15").Value = -199718
$ws.Range("N15").Value = -70583

$ws.Range("H22").Value = 2993.3333
$ws.Range("J22").Value = 2993.3333
$ws.Range("L22").Value = 2993.3333
$ws.Range("N22").Value = -3579.3333
